# The commit "Deploy the implementation guide" refreshes the generated
# ValueSet metadata sheet: the resource's Status flips from "active" to
# "draft" and the Date stamp moves to the new publish time. (The diff's
# styles.xml hunk just records cellXfs gaining the applyAlignment="true"
# flag for the alignment that was already declared on those two styles -
# i.e. the already-authored vertical/wrap alignment becoming effective;
# there is no visible formatting left to author here.)

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2023-08-01T16:12:28+00:00"
